$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Resize / reposition the red "Rectangulo 1" text-box anchor (the one that
#    wraps "{month} - {year}") - it got a bit taller (36pt -> 41.25pt) and its
#    vertical offset was re-derived by Word by a hair (454025 -> 454024 EMU).
# ---------------------------------------------------------------------------
$shape = $d.Shapes.Item(1)
$shape.Height = 41.25
$shape.Top = 35.74992125984252

# ---------------------------------------------------------------------------
# 2. Split "({%image}" into "({%" / "image}" with a gramStart/gramEnd proof
#    error pair in between (exactly what Word's grammar checker inserts when
#    it flags the run boundary).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("({%image}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $rng.Delete()
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
      '<w:body>' +
      '<w:p w14:paraId="43D70928" w14:textId="37820ADA" w:rsidR="00715780" w:rsidRPr="0021503C" w:rsidRDefault="0021503C" w:rsidP="00715780">' +
      '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>({%</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>image}</w:t></w:r>' +
      '</w:p>' +
      '</w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3. Split "{LeadAuditor}" into "{" / "LeadAuditor" / "}" with a
#    spellStart/spellEnd proof error pair bracketing the middle word - the
#    surrounding runs ("Leader Auditor: " and ".") are left untouched.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("{LeadAuditor}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $rng2.Delete()
    $xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
      '<w:body>' +
      '<w:p w14:paraId="6FF9583A" w14:textId="50C6ECB4" w:rsidR="006166E3" w:rsidRPr="006166E3" w:rsidRDefault="006166E3" w:rsidP="006166E3">' +
      '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
      '<w:r w:rsidRPr="006166E3"><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Leader Auditor: </w:t></w:r>' +
      '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>LeadAuditor</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r>' +
      '<w:r w:rsidRPr="006166E3"><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r>' +
      '</w:p>' +
      '</w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'
    $rng2.InsertXML($xml2)
}
